$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("location")

# A3 was its own unique string "Barracks " (trailing space); dedupe it onto the
# shared "Barracks" string used elsewhere (C5/B6/E7), which drops the now-unused
# shared-string entry.
$ws.Range("A3").Value = "Barracks"

# Row 7: "Range" moves from D7 to C7.
$ws.Range("C7").Value = $ws.Range("D7").Value2
$ws.Range("D7").Value = ""

# Row 8: "Courtyard" moves from E8 to B8.
$ws.Range("B8").Value = $ws.Range("E8").Value2
$ws.Range("E8").Value = ""

# Active-cell selection on the location sheet moves from F2 to A3.
$ws.Range("A3").Select() | Out-Null
